$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ E=2; F=0.6666666666666666; G=0.638706; H=1.916118; I=0.5586654432763536; J=0.5586654432763536; M=0.7109393333333333; N=2.132818; O=0.002867418779909113; P=0.002875587087952277; Q=0.454081217836; R=4.086730960524; S=0.001601927783736865; T=0.001606491135170618 }
    3 = @{ E=2; F=0.6666666666666666; G=0.638706; H=1.916118; I=0.5586654432763536; J=0.5586654432763536; M=132.4704766666667; N=397.4114300000001; O=0.5342907823042267; P=0.5358127963626762; Q=84.60968826986002; R=761.4871944287401; S=0.2984897967344605; T=0.2993400933930971 }
    4 = @{ E=2; F=0.6666666666666666; G=0.638706; H=1.916118; I=0.5586654432763536; J=0.5586654432763536; K=2; L=1; M=2.112848; N=4.225696; O=0.008521711699207447; P=0.005697324785899025; Q=1.349488694688; R=8.096932168127999; S=0.004760785843911017; T=0.003182898477003635 }
    5 = @{ E=2; F=0.6666666666666666; G=0.638706; H=1.916118; I=0.5586654432763536; J=0.5586654432763536; M=112.642779; N=337.928337; O=0.4543200872166568; P=0.4556142917634725; Q=71.945618803974; R=647.510569235766; S=0.2538129329142452; T=0.2545359602710822 }
    6 = @{ G=0.5045649999999999; H=1.513695; I=0.4413345567236464; J=0.4413345567236464; M=0.7109393333333333; N=2.132818; O=0.002867418779909113; P=0.002875587087952277; Q=0.3587151047233332; R=3.22843594251; S=0.001265490996172247; T=0.001269095952781659 }
    7 = @{ G=0.5045649999999999; H=1.513695; I=0.4413345567236464; J=0.4413345567236464; M=132.4704766666667; N=397.4114300000001; O=0.5342907823042267; P=0.5358127963626762; Q=66.83996605931667; R=601.55969453385; S=0.2358009855697661; T=0.2364727029695791 }
    8 = @{ G=0.5045649999999999; H=1.513695; I=0.4413345567236464; J=0.4413345567236464; K=2; L=1; M=2.112848; N=4.225696; O=0.008521711699207447; P=0.005697324785899025; Q=1.06606915112; R=6.39641490672; S=0.00376092585529643; T=0.00251442630889539 }
    9 = @{ G=0.5045649999999999; H=1.513695; I=0.4413345567236464; J=0.4413345567236464; M=112.642779; N=337.928337; O=0.4543200872166568; P=0.4556142917634725; Q=56.835603786135; R=511.520434075215; S=0.2005071543024116; T=0.2010783314923902 }
}

foreach ($rowNum in $data.Keys) {
    $rowData = $data[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
